$wb = $excel.ActiveWorkbook

# ---- Sheet "routes" ----
$routes = $wb.Worksheets.Item("routes")

# Row 2: /user /login - current task updated
$routes.Range("F2").Value = "current task"
$routes.Range("G2").Value = 45556
$routes.Range("H2").Value = "to check option of cookies and logout"

# Row 3: /register - remarks updated
$routes.Range("D3").Value = "register a new user to course"

# Row 5: /add (faculty) - status moved to done, remarks + endpoint + test suite added
$routes.Range("F5").Value = "done"
$routes.Range("G5").Value = 45556
$routes.Range("H5").Value = "admin or super admin can create a new faculty (dependencies: valid departmentId and subjectId - to be tested)"
$routes.Range("I5").Value = "/api/user/add"
$routes.Range("J5").Value = "user-add.test.js"

# Row 6: /find -> /get
$routes.Range("B6").Value = "/get"

# ---- Sheet "tasks" ----
$tasks = $wb.Worksheets.Item("tasks")

# Row 5 remarks updated
$tasks.Range("D5").Value = "Tested with jest and supertest. api doc to be created."

# New row 6 - add faculty task completed
$tasks.Range("A6").Value = 45556
$tasks.Range("B6").Value = "Jasdeep"
$tasks.Range("C6").Value = "Wrote and tested /api/user/add route to add faculty"
$tasks.Range("D6").Value = "stand alone testing done. Yet to perform integration with department and subject models."

# Match the date cell style used by the rest of column A (numFmtId 16 date format)
$tasks.Range("A6").NumberFormat = $tasks.Range("A5").NumberFormat

[void]$tasks.Activate()
$tasks.Range("A6").Select() | Out-Null

# Re-activate "routes" as the active/selected tab, with H5 selected
[void]$routes.Activate()
$routes.Range("H5").Select() | Out-Null
